# Add two new match rows (45 and 46) to the Premijer Liga BiH 2023-2024 sheet,
# matching the data appended upstream (Tuzla City vs Zrinjski, and
# Siroki Brijeg vs Borac Banja Luka).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 45 : Tuzla City 2 - 4 Zrinjski ----
$r = 45
$ws.Cells.Item($r, 1).Value = 44
$ws.Cells.Item($r, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($r, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45196.625
$ws.Cells.Item($r, 6).Value = "Tuzla City"
$ws.Cells.Item($r, 7).Value = 2
$ws.Cells.Item($r, 8).Value = "Zrinjski"
$ws.Cells.Item($r, 9).Value = 4
$ws.Cells.Item($r, 10).Value = 3.53
$ws.Cells.Item($r, 11).Value = "04/08/2023 04:12"
$ws.Cells.Item($r, 12).Value = 5.98
$ws.Cells.Item($r, 13).Value = "27/09/2023 14:59"
$ws.Cells.Item($r, 14).Value = 3.17
$ws.Cells.Item($r, 15).Value = "04/08/2023 04:12"
$ws.Cells.Item($r, 16).Value = 4.11
$ws.Cells.Item($r, 17).Value = "27/09/2023 14:59"
$ws.Cells.Item($r, 18).Value = 1.97
$ws.Cells.Item($r, 19).Value = "04/08/2023 04:12"
$ws.Cells.Item($r, 20).Value = 1.51
$ws.Cells.Item($r, 21).Value = "27/09/2023 14:59"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/tuzla-city-zrinjski/QTjlSnln/"

# ---- Row 46 : Siroki Brijeg 0 - 2 Borac Banja Luka ----
$r = 46
$ws.Cells.Item($r, 1).Value = 45
$ws.Cells.Item($r, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($r, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45196.77083333334
$ws.Cells.Item($r, 6).Value = "Siroki Brijeg"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = "Borac Banja Luka"
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 2.73
$ws.Cells.Item($r, 11).Value = "26/09/2023 20:42"
$ws.Cells.Item($r, 12).Value = 3.15
$ws.Cells.Item($r, 13).Value = "27/09/2023 18:16"
$ws.Cells.Item($r, 14).Value = 3.2
$ws.Cells.Item($r, 15).Value = "26/09/2023 20:42"
$ws.Cells.Item($r, 16).Value = 2.63
$ws.Cells.Item($r, 17).Value = "27/09/2023 18:16"
$ws.Cells.Item($r, 18).Value = 2.5
$ws.Cells.Item($r, 19).Value = "26/09/2023 20:42"
$ws.Cells.Item($r, 20).Value = 2.68
$ws.Cells.Item($r, 21).Value = "27/09/2023 18:16"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/siroki-brijeg-borac-banja-luka/6Nzl3FYa/"

# Replicate the formatting of the previous data row (44) for the styled
# columns: A ("Indice" - bold/centered/bordered) and E ("data_partida" -
# date/time number format).
$ws.Range("A44").Copy()
$ws.Range("A45:A46").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E44").Copy()
$ws.Range("E45:E46").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
